$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 67; $r++) {
    $ws.Cells.Item($r, 3).Value = 7590
}

for ($r = 68; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7573
}
